# The experiment workflow was fixed/re-run: for each results sheet (NBR and
# BAR), the first 4 cutoff steps of the scan were dropped. The "Cutoff"
# column (A) keeps its original sequential numbering (0..14), while the
# "Cutoff step"/"Reaction_number" columns (B, C) now show what used to be
# rows 6-20, and the table shrinks from 19 to 15 data rows.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Grab the old B6:C20 values (the ones to keep) before overwriting anything.
    $tail = $ws.Range("B6:C20").Value()

    # Shift them up into B2:C16, overwriting the first 15 rows of B/C data.
    $ws.Range("B2:C16").Value = $tail

    # Drop the now-duplicated trailing rows (17-20), shrinking the table
    # and the sheet dimension down to A1:C16.
    $ws.Range("A17:C20").EntireRow.Delete()
}
